$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (names, links, percentage strings, and non-numeric-looking prices)
$textUpdates = @{
    'D2' = '29.286.45'
    'E2' = '  -0.14%  '
    'D3' = '1.867.59'
    'E3' = '  +0.36%  '
    'E4' = '  -0.02%  '
    'E5' = '  +2.93%  '
    'E6' = '  +1.29%  '
    'E7' = '  +0.04%  '
    'E8' = '  -0.05%  '
    'E9' = '  +1.30%  '
    'E10' = '  +2.20%  '
    'E11' = '  +1.09%  '
    'E12' = '  +1.12%  '
    'E13' = '  +0.41%  '
    'D14' = '1.814.40'
    'E14' = '  -2.14%  '
    'E15' = '  +1.75%  '
    'D16' = '29.277.15'
    'E16' = '  -0.25%  '
    'E17' = '  +1.10%  '
    'E18' = '  +0.88%  '
    'E19' = '  +0.45%  '
    'E20' = '  +0.43%  '
    'D21' = '2.108.15'
    'E21' = '  -0.31%  '
    'E22' = '  -0.06%  '
    'E23' = '  +6.51%  '
    'E24' = '  +0.02%  '
    'E25' = '  +11.37%  '
    'E26' = '  -0.31%  '
    'E27' = '  +0.76%  '
    'E28' = '  +0.92%  '
    'E29' = '  -1.92%  '
    'E30' = '  +1.67%  '
    'E31' = '  +2.45%  '
    'E32' = '  +1.76%  '
    'E33' = '  +0.82%  '
    'E34' = '  +1.78%  '
    'E35' = '  +0.57%  '
    'E36' = '  +3.22%  '
    'E37' = '  +0.31%  '
    'E38' = '  +0.81%  '
    'E39' = '  +0.47%  '
    'D40' = '1.170.86'
    'E40' = '  -0.10%  '
    'E41' = '  -1.56%  '
    'E42' = '  +1.40%  '
    'E43' = '  +0.94%  '
    'E44' = '  +0.08%  '
    'E45' = '  +0.03%  '
    'B46' = 'Mantle'
    'C46' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'E46' = '  -1.18%  '
    'B47' = 'RocketPoolETH'
    'C47' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D47' = '1.998.39'
    'E47' = '  -0.71%  '
    'E48' = '  +1.86%  '
    'B49' = 'BabyDogeCoin'
    'C49' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'E49' = '  +1.10%  '
    'B50' = 'SynthetixNetwork'
    'C50' = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
    'E50' = '  +6.52%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E51' = '  +1.64%  '
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price updates that look like plain numbers - must be forced to remain text
# so Excel does not coerce them into numeric values (losing formatting like trailing zeros).
$numericTextUpdates = @{
    'D5' = '0.7231'
    'D6' = '240.95'
    'D8' = '0.07839'
    'D9' = '0.3091'
    'D11' = '0.08246'
    'D12' = '0.7227'
    'D13' = '5.237'
    'D15' = '90.83'
    'D17' = '5.861'
    'D18' = '243.90'
    'D19' = '0.000007813'
    'D20' = '13.22'
    'D22' = '1.000'
    'D23' = '8.002'
    'D24' = '1.001'
    'D25' = '0.1590'
    'D26' = '162.17'
    'D27' = '8.967'
    'D29' = '1.347'
    'D31' = '4.407'
    'D32' = '4.106'
    'D33' = '0.05208'
    'D34' = '1.937'
    'D35' = '1.187'
    'D36' = '0.7286'
    'D37' = '2.682'
    'D39' = '2.703'
    'D41' = '0.9039'
    'D42' = '6.096'
    'D43' = '72.12'
    'D44' = '1.001'
    'D45' = '101.69'
    'D46' = '0.5287'
    'D48' = '1.781'
    'D49' = '0.00000000120'
    'D50' = '2.893'
    'D51' = '9.299'
}
foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$ref]
    $cell.Style = $origStyle
}
